$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) sometimes holds numeric-looking text (e.g. "0.160", "588.61").
# Typing such text into a cell via .Value lets Excel auto-convert it into a real
# number, which would silently drop significant trailing zeros (e.g. "0.160" ->
# 0.16) or switch to scientific notation. Mark those specific cells as Text first
# so the exact original price string is preserved, same as the scraped source data.
$priceTextCells = @(
    "D4", "D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D17", "D19", "D20", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D49", "D50"
)
foreach ($cellRef in $priceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.251.00"
$ws.Range("E2").Value = "  -3.59%  "
$ws.Range("D3").Value = "3.651.04"
$ws.Range("E3").Value = "  -5.14%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "588.61"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").Value = "178.56"
$ws.Range("E6").Value = "  +5.45%  "
$ws.Range("D7").Value = "3.647.53"
$ws.Range("E7").Value = "  -5.10%  "
$ws.Range("D8").Value = "0.628"
$ws.Range("E8").Value = "  -5.56%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "0.711"
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("D11").Value = "0.160"
$ws.Range("E11").Value = "  -8.87%  "
$ws.Range("D12").Value = "55.36"
$ws.Range("E12").Value = "  +4.46%  "
$ws.Range("D13").Value = "0.0000290"
$ws.Range("E13").Value = "  -9.44%  "
$ws.Range("D14").Value = "10.63"
$ws.Range("E14").Value = "  -5.36%  "
$ws.Range("D15").Value = "4.222.34"
$ws.Range("E15").Value = "  -5.21%  "
$ws.Range("D16").Value = "3.644.97"
$ws.Range("E16").Value = "  -5.43%  "
$ws.Range("D17").Value = "19.25"
$ws.Range("E17").Value = "  -9.11%  "
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").Value = "1.12"
$ws.Range("E19").Value = "  -6.86%  "
$ws.Range("D20").Value = "12.69"
$ws.Range("E20").Value = "  -8.20%  "
$ws.Range("D21").Value = "67.927.97"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").Value = "408.65"
$ws.Range("E22").Value = "  -6.35%  "
$ws.Range("D23").Value = "4.56"
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("D24").Value = "88.15"
$ws.Range("E24").Value = "  -6.30%  "
$ws.Range("E25").Value = "  -8.82%  "
$ws.Range("D26").Value = "12.65"
$ws.Range("E26").Value = "  -8.37%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "10.77"
$ws.Range("E27").Value = "  -6.61%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "3.87"
$ws.Range("E28").Value = "  -2.80%  "
$ws.Range("D29").Value = "6.06"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "9.46"
$ws.Range("E30").Value = "  -10.08%  "
$ws.Range("D31").Value = "32.47"
$ws.Range("E31").Value = "  -6.96%  "
$ws.Range("D32").Value = "7.17"
$ws.Range("E32").Value = "  -12.80%  "
$ws.Range("D33").Value = "12.28"
$ws.Range("E33").Value = "  -8.80%  "
$ws.Range("D34").Value = "0.117"
$ws.Range("E34").Value = "  -6.90%  "
$ws.Range("D35").Value = "64.51"
$ws.Range("E35").Value = "  -6.23%  "
$ws.Range("D36").Value = "600.52"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("D37").Value = "42.70"
$ws.Range("E37").Value = "  -10.92%  "
$ws.Range("D38").Value = "0.0₃0882"
$ws.Range("E38").Value = "  -9.60%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "0.395"
$ws.Range("E40").Value = "  -8.75%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  -6.77%  "
$ws.Range("D43").Value = "3.01"
$ws.Range("E43").Value = "  -6.57%  "
$ws.Range("D44").Value = "2.69"
$ws.Range("E44").Value = "  -6.73%  "
$ws.Range("D45").Value = "0.0437"
$ws.Range("E45").Value = "  -6.75%  "
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  -11.66%  "
$ws.Range("D47").Value = "0.134"
$ws.Range("E47").Value = "  -6.22%  "
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("D49").Value = "8.94"
$ws.Range("E49").Value = "  -9.96%  "
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("D51").Value = "2.698.05"
$ws.Range("E51").Value = "  -7.01%  "
